# Add the "ODI Bowling Extra" sheet (sheetId 5) after the last existing sheet
# ("ODI Batting Extra"), matching the ODI-extras layout already used by the
# workbook (MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL).

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# --- Header row -------------------------------------------------------
$headers = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
}

# --- Data rows ----------------------------------------------------------
# columns: row, MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL ($null = blank)
$data = @(
    ,@(2,  "4326", $null, $null)
    ,@(3,  "4331", "0",   "10.00%")
    ,@(4,  "4336", "0",   "20.00%")
    ,@(5,  "4342", "3",   "20.00%")
    ,@(6,  "4346", $null, $null)
    ,@(7,  "4354", "0",   "30.00%")
    ,@(8,  "4355", $null, $null)
    ,@(9,  "4401", $null, $null)
    ,@(10, "4405", $null, $null)
    ,@(11, "4429", "0",   "10.00%")
    ,@(12, "4430", "1",   "30.00%")
    ,@(13, "4431", "0",   "20.00%")
    ,@(14, "4469", "5",   "40.00%")
    ,@(15, "4471", $null, $null)
    ,@(16, "4663", $null, $null)
    ,@(17, "4666", $null, $null)
    ,@(18, "4699", "0",   $null)
    ,@(19, "4700", "1",   "10.00%")
    ,@(20, "4711", "0",   "10.00%")
    ,@(21, "4717", "0",   "10.00%")
)

foreach ($row in $data) {
    $r = $row[0]
    $matchCode = $row[1]
    $maidenOvers = $row[2]
    $percentWickets = $row[3]

    # Force text storage (these look numeric/percent-like but must stay
    # strings, matching the sibling "ODI Batting Extra" sheet's layout).
    $aCell = $ws.Cells.Item($r, 1)
    $aCell.NumberFormat = "@"
    $aCell.Value = $matchCode

    if ($maidenOvers -ne $null) {
        $bCell = $ws.Cells.Item($r, 2)
        $bCell.NumberFormat = "@"
        $bCell.Value = $maidenOvers
    }

    if ($percentWickets -ne $null) {
        $cCell = $ws.Cells.Item($r, 3)
        $cCell.NumberFormat = "@"
        $cCell.Value = $percentWickets
    }
}

[void]$ws.Range("A1").Select()

Write-Output "Added 'ODI Bowling Extra' sheet with header + 20 data rows"
